$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# C10 ("Integer min" for rule R20) changes from 18 to 1
$ws.Range("C10").Value = 1
